# "drop down stat-list added"
# Add three new response rows (11-13) to the "Ответы" sheet, all referring to
# the new ПВИ location "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ответы")

$uniqueKey = 99295
$piiLocation = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)"
$answerDate = 44165

$newRows = @(
    @{ Row = 11; Idx = 9;  Question = "Довольны ли Вы качеством предоставляемого питания?" },
    @{ Row = 12; Idx = 10; Question = "Устраивают ли Вас бытовые условия? (питьевой режим, температура в помещении, досуг и психологический климат)" },
    @{ Row = 13; Idx = 11; Question = "Довольны ли Вы работой обслуживающего персонала?" }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.Idx
    $ws.Cells.Item($r, 2).Value = $uniqueKey
    $ws.Cells.Item($r, 3).Value = $item.Question
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = "None"
    $ws.Cells.Item($r, 6).Value = $answerDate
    $ws.Cells.Item($r, 7).Value = $piiLocation

    # Copy formatting from the last pre-existing data row (10) so the new
    # rows match the sheet's existing style (index style for col A, date
    # format for col F) without introducing any new cell-style entries.
    $ws.Cells.Item(10, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item(10, 6).Copy() | Out-Null
    $ws.Cells.Item($r, 6).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false
